# Updated cryptos list refresh (prices + 1h volume deltas), mirroring the
# "Updated cryptos list ... with GitHub Actions" scheduled commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel keeps trailing
# zeros / precision as literal text instead of coercing to a Double.
$textCells = @("D4","D5","D7","D9","D11","D12","D13","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "26.341.42"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.713.50"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").Value = "0.9960"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "240.14"
$ws.Range("E5").Value = "  -2.66%  "
$ws.Range("D7").Value = "0.4871"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").Value = "0.06169"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "1.712.16"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "0.06949"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "15.51"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "0.5975"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "4.462"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "76.46"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "0.9964"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "26.237.89"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "0.9960"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "0.000007090"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "1.932.33"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "4.399"
$ws.Range("E22").Value = "  -4.13%  "
$ws.Range("D23").Value = "8.423"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").Value = "5.042"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").Value = "136.08"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").Value = "15.18"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "1.393"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Value = "1.728"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "105.32"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("D31").Value = "0.07932"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "3.605"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "0.04445"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.595"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "0.9900"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6170"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9379"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "1.990"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.367"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9957"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01471"
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "99.68"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.357"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3799"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.812"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1150"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05348"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "30.53"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.688"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "51.18"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.208"
$ws.Range("E51").Value = "  -4.80%  "
